$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row at the top, pushing all existing rows down by one.
$ws.Rows("1:1").Insert()

# Populate the new header row with the category labels.
$ws.Range("A1").Value = "category"
$ws.Range("B1").Value = "all"
$ws.Range("C1").Value = "smi"
$ws.Range("D1").Value = "men"
$ws.Range("E1").Value = "women"
